# "params" sheet gains a new "description" header column (M1), right after
# the existing "ui variable" column (L1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

$ws.Range("M1").Value = "description"

# Leave the freshly added header cell selected, matching the edited
# selection/activeCell state of the source file.
$ws.Range("M1").Select() | Out-Null
